# Updated cryptos list - apply latest price/volume snapshot values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep Price column (D) as text so values like "46.404.19" or "308.54"
# are not re-interpreted as numbers/dates when the value is assigned.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2 - Bitcoin
$ws.Range("D2").Value = "46.404.19"
$ws.Range("E2").Value = "  +1.40%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.613.65"
$ws.Range("E3").Value = "  +3.71%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.01%  "

# Row 5 - BNB
$ws.Range("D5").Value = "308.54"
$ws.Range("E5").Value = "  +3.12%  "

# Row 6 - Solana
$ws.Range("D6").Value = "100.45"
$ws.Range("E6").Value = "  +3.22%  "

# Row 7 - XRP
$ws.Range("E7").Value = "  +2.44%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  -0.11%  "

# Row 9 - Cardano
$ws.Range("D9").Value = "0.581"
$ws.Range("E9").Value = "  +6.66%  "

# Row 10 - Avalanche
$ws.Range("D10").Value = "39.61"
$ws.Range("E10").Value = "  +7.24%  "

# Row 11 - Dogecoin
$ws.Range("E11").Value = "  +5.50%  "

# Row 12 - OKB
$ws.Range("D12").Value = "54.28"
$ws.Range("E12").Value = "  +0.50%  "

# Row 13 - Polkadot
$ws.Range("E13").Value = "  +6.67%  "

# Row 14 - WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "3.009.68"
$ws.Range("E14").Value = "  +3.44%  "

# Row 15 - TRON
$ws.Range("E15").Value = "  +1.50%  "

# Row 16 - WrappedEther
$ws.Range("D16").Value = "2.608.52"
$ws.Range("E16").Value = "  +3.29%  "

# Row 17 - Polygon
$ws.Range("D17").Value = "0.923"
$ws.Range("E17").Value = "  +5.07%  "

# Row 18 - Chainlink
$ws.Range("D18").Value = "14.96"
$ws.Range("E18").Value = "  +2.78%  "

# Row 19 - WrappedBTC
$ws.Range("D19").Value = "46.548.94"
$ws.Range("E19").Value = "  +1.45%  "

# Row 20 - ShibaInu
$ws.Range("E20").Value = "  +4.80%  "

# Row 21 - InternetComputer(DFINITY)
$ws.Range("D21").Value = "12.97"
$ws.Range("E21").Value = "  -2.48%  "

# Row 22 - Uniswap
$ws.Range("E22").Value = "  +3.82%  "

# Row 23 - Litecoin
$ws.Range("D23").Value = "71.59"
$ws.Range("E23").Value = "  +4.11%  "

# Row 24 - BitcoinCash
$ws.Range("D24").Value = "274.52"
$ws.Range("E24").Value = "  +9.96%  "

# Row 25 - PancakeSwap
$ws.Range("E25").Value = "  +6.12%  "

# Row 26 - ImmutableX
$ws.Range("D26").Value = "2.16"
$ws.Range("E26").Value = "  +5.61%  "

# Row 27 - EthereumClassic
$ws.Range("D27").Value = "28.96"
$ws.Range("E27").Value = "  +27.72%  "

# Row 28 - Dai
$ws.Range("E28").Value = "  -0.04%  "

# Row 29 - LEO
$ws.Range("E29").Value = "  -1.09%  "

# Row 30 - Cosmos
$ws.Range("E30").Value = "  +5.39%  "

# Rows 31/32 - InjectiveProtocol and Toncoin swap places (ranking order changed)
$ws.Range("B31").Value = "InjectiveProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D31").Value = "39.14"
$ws.Range("E31").Value = "  -4.15%  "

$ws.Range("B32").Value = "Toncoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D32").Value = "2.23"
$ws.Range("E32").Value = "  -0.47%  "

# Row 33 - Filecoin
$ws.Range("E33").Value = "  +10.21%  "

# Row 34 - LidoDAOToken
$ws.Range("E34").Value = "  -5.33%  "

# Row 35 - ARBITRUM
$ws.Range("E35").Value = "  +2.74%  "

# Row 36 - WEMIXToken
$ws.Range("E36").Value = "  +2.52%  "

# Row 37 - Hedera
$ws.Range("E37").Value = "  +4.60%  "

# Row 38 - Monero
$ws.Range("D38").Value = "151.10"
$ws.Range("E38").Value = "  +1.29%  "

# Row 39 - Kaspa
$ws.Range("E39").Value = "  +4.95%  "

# Row 40 - Stellar
$ws.Range("E40").Value = "  +4.71%  "

# Row 41 - EnergySwap
$ws.Range("D41").Value = "23.38"
$ws.Range("E41").Value = "  +40.42%  "

# Row 42 - Celestia
$ws.Range("D42").Value = "15.92"
$ws.Range("E42").Value = "  +0.94%  "

# Row 43 - NEARProtocol
$ws.Range("E43").Value = "  +8.45%  "

# Row 44 - VeChain
$ws.Range("E44").Value = "  +6.83%  "

# Row 45 - RenderToken
$ws.Range("D45").Value = "4.11"
$ws.Range("E45").Value = "  +0.07%  "

# Row 46 - Maker
$ws.Range("D46").Value = "2.128.99"
$ws.Range("E46").Value = "  +5.64%  "

# Row 47 - FirstDigitalUSD
$ws.Range("D47").Value = "0.997"
$ws.Range("E47").Value = "  -0.09%  "

# Row 48 - BitcoinSV
$ws.Range("D48").Value = "93.70"
$ws.Range("E48").Value = "  +2.69%  "

# Row 49 - FraxShare
$ws.Range("E49").Value = "  +8.43%  "

# Row 50 - Stacks
$ws.Range("D50").Value = "1.79"
$ws.Range("E50").Value = "  -0.72%  "

# Row 51 - Aave
$ws.Range("D51").Value = "109.26"
$ws.Range("E51").Value = "  +3.12%  "
